$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.146793842315674
$ws.Range("B1").Value = 1.381637692451477
$ws.Range("C1").Value = 1.853913068771362
$ws.Range("D1").Value = 3.479630470275879
$ws.Range("E1").Value = 1.849699378013611
